# Generate Report for Handoff
# Adds two new rows (a dependency .png and its .md "parent" + another
# dependency .png) to each of the three sheets, updates the existing
# handoff row's timestamp/filename, and wires up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Existing row 2 gets a new source-file name + refreshed handoff date.
$ov.Range("A2").Value = "07978c49-e6f5-4302-ab05-7e860d13dbb2.png"
$ov.Range("D2").Value = "2016-03-22 15:07:35"
$ov.Hyperlinks.Item(1).Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/07978c49-e6f5-4302-ab05-7e860d13dbb2.png", "", "", "07978c49-e6f5-4302-ab05-7e860d13dbb2.png") | Out-Null

# New row 3: the markdown file that references the two pngs.
$ov.Range("A3").Value = "0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-22 15:07:35"
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/0ad07800-945f-4fbe-a25b-fc6a016db211.md", "", "", "0ad07800-945f-4fbe-a25b-fc6a016db211.md") | Out-Null

# New row 4: the second png (handed off alongside the markdown).
$ov.Range("A4").Value = "c19546b8-c530-4464-a08b-e860884a8d63.png"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-22 15:07:35"
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/c19546b8-c530-4464-a08b-e860884a8d63.png", "", "", "c19546b8-c530-4464-a08b-e860884a8d63.png") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "07978c49-e6f5-4302-ab05-7e860d13dbb2.png"
$zh.Range("B2").Value = ".png"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "0166c388fc6a6bf4127c6577e8bc88b74670fddb.png"
$zh.Range("E2").Value = "2016-03-22 15:07:32"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "IsDependency"
$zh.Range("K2").Value = "e2e\0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$zh.Hyperlinks.Item(2).Delete()
$zh.Hyperlinks.Item(1).Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/07978c49-e6f5-4302-ab05-7e860d13dbb2.png", "", "", "07978c49-e6f5-4302-ab05-7e860d13dbb2.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/006331faecbd2b4d6bc40e4470654897db1d5160/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0166c388fc6a6bf4127c6577e8bc88b74670fddb.png", "", "", "0166c388fc6a6bf4127c6577e8bc88b74670fddb.png") | Out-Null

$zh.Range("A3").Value = "0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-22 15:07:32"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/0ad07800-945f-4fbe-a25b-fc6a016db211.md", "", "", "0ad07800-945f-4fbe-a25b-fc6a016db211.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/006331faecbd2b4d6bc40e4470654897db1d5160/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.zh-cn.xlf", "", "", "0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.zh-cn.xlf") | Out-Null

$zh.Range("A4").Value = "c19546b8-c530-4464-a08b-e860884a8d63.png"
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "4bf66b52de4364a5900704873331e1d51a7831d7.png"
$zh.Range("E4").Value = "2016-03-22 15:07:32"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("J4").Value = "IsDependency"
$zh.Range("K4").Value = "e2e\0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/c19546b8-c530-4464-a08b-e860884a8d63.png", "", "", "c19546b8-c530-4464-a08b-e860884a8d63.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/006331faecbd2b4d6bc40e4470654897db1d5160/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4bf66b52de4364a5900704873331e1d51a7831d7.png", "", "", "4bf66b52de4364a5900704873331e1d51a7831d7.png") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "07978c49-e6f5-4302-ab05-7e860d13dbb2.png"
$de.Range("B2").Value = ".png"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "0166c388fc6a6bf4127c6577e8bc88b74670fddb.png"
$de.Range("E2").Value = "2016-03-22 15:07:35"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "IsDependency"
$de.Range("K2").Value = "e2e\0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$de.Hyperlinks.Item(2).Delete()
$de.Hyperlinks.Item(1).Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/07978c49-e6f5-4302-ab05-7e860d13dbb2.png", "", "", "07978c49-e6f5-4302-ab05-7e860d13dbb2.png") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9da81897d1daff9579ecc03cc75091c5798701c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0166c388fc6a6bf4127c6577e8bc88b74670fddb.png", "", "", "0166c388fc6a6bf4127c6577e8bc88b74670fddb.png") | Out-Null

$de.Range("A3").Value = "0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.de-de.xlf"
$de.Range("E3").Value = "2016-03-22 15:07:35"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/0ad07800-945f-4fbe-a25b-fc6a016db211.md", "", "", "0ad07800-945f-4fbe-a25b-fc6a016db211.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9da81897d1daff9579ecc03cc75091c5798701c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.de-de.xlf", "", "", "0ad07800-945f-4fbe-a25b-fc6a016db211.a803b81c6c20c4a32c67038bb8069d88bb5df34e.de-de.xlf") | Out-Null

$de.Range("A4").Value = "c19546b8-c530-4464-a08b-e860884a8d63.png"
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "4bf66b52de4364a5900704873331e1d51a7831d7.png"
$de.Range("E4").Value = "2016-03-22 15:07:35"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("J4").Value = "IsDependency"
$de.Range("K4").Value = "e2e\0ad07800-945f-4fbe-a25b-fc6a016db211.md"
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/379022beba32f37b4471bffdb419ceb631de7d55/e2e/c19546b8-c530-4464-a08b-e860884a8d63.png", "", "", "c19546b8-c530-4464-a08b-e860884a8d63.png") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9da81897d1daff9579ecc03cc75091c5798701c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4bf66b52de4364a5900704873331e1d51a7831d7.png", "", "", "4bf66b52de4364a5900704873331e1d51a7831d7.png") | Out-Null

Write-Host "Report generated for handoff"
